$wb = $excel.ActiveWorkbook

# --- RUNMANAGER sheet (small "execute" table) ---
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws1.Range("C2").Value = "yes"
$ws1.Range("E2").Value = "'2"
$ws1.Range("E4").Select()

# --- DATA sheet (the "run manager" rows table) ---
$ws2 = $wb.Worksheets.Item("DATA")
$ws2.Range("B4").Value = "no"
$ws2.Range("B5").Value = "no"
$ws2.Range("B7").Value = "yes"
$ws2.Range("B9").Value = "yes"
$ws2.Range("A9").Select()
